$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price (D) and Volume (E) columns hold numeric-looking text
# (e.g. "26.170.57", "4.01") that must stay text, not become numbers.
# Temporarily force the "Text" number format on the data range so that
# assigning these values keeps them as strings, then restore the
# original (default/"Normal") style so the saved styles match the
# original workbook's cell formatting.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '26.170.57'
$ws.Range("E2").Value = '  +3.58%  '
$ws.Range("D3").Value = '1.603.85'
$ws.Range("E3").Value = '  +3.40%  '
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("D5").Value = '212.88'
$ws.Range("E5").Value = '  +2.97%  '
$ws.Range("E6").Value = '  -0.24%  '
$ws.Range("E7").Value = '  +1.75%  '
$ws.Range("E8").Value = '  +2.54%  '
$ws.Range("D9").Value = '0.0616'
$ws.Range("E9").Value = '  +0.89%  '
$ws.Range("D10").Value = '18.03'
$ws.Range("E10").Value = '  +1.45%  '
$ws.Range("D11").Value = '0.0818'
$ws.Range("E11").Value = '  +4.88%  '
$ws.Range("D12").Value = '1.825.33'
$ws.Range("E12").Value = '  +3.39%  '
$ws.Range("D13").Value = '1.600.02'
$ws.Range("E13").Value = '  +3.01%  '
$ws.Range("D14").Value = '4.01'
$ws.Range("E14").Value = '  +0.49%  '
$ws.Range("D15").Value = '0.512'
$ws.Range("E15").Value = '  +1.40%  '
$ws.Range("D16").Value = '26.145.76'
$ws.Range("E16").Value = '  +3.65%  '
$ws.Range("D17").Value = '60.52'
$ws.Range("E17").Value = '  +3.16%  '
$ws.Range("D18").Value = '0.0₃0723'
$ws.Range("E18").Value = '  +2.14%  '
$ws.Range("E19").Value = '  -0.25%  '
$ws.Range("D20").Value = '204.30'
$ws.Range("E20").Value = '  +9.84%  '
$ws.Range("D21").Value = '4.24'
$ws.Range("E21").Value = '  +3.35%  '
$ws.Range("D22").Value = '9.31'
$ws.Range("E22").Value = '  +0.63%  '
$ws.Range("D23").Value = '5.99'
$ws.Range("E23").Value = '  +2.87%  '
$ws.Range("D24").Value = '1.84'
$ws.Range("E24").Value = '  +12.52%  '
$ws.Range("D25").Value = '141.89'
$ws.Range("E25").Value = '  +1.84%  '
$ws.Range("E26").Value = '  -0.29%  '
$ws.Range("E27").Value = '  -4.42%  '
$ws.Range("D28").Value = '15.21'
$ws.Range("E28").Value = '  +2.48%  '
$ws.Range("D29").Value = '6.43'
$ws.Range("E29").Value = '  +0.76%  '
$ws.Range("D31").Value = '0.0472'
$ws.Range("E31").Value = '  +1.77%  '
$ws.Range("D32").Value = '3.12'
$ws.Range("E32").Value = '  +3.04%  '
$ws.Range("D33").Value = '2.99'
$ws.Range("E33").Value = '  +0.87%  '
$ws.Range("E34").Value = '  +1.31%  '
$ws.Range("E35").Value = '  +1.60%  '

# Rows 36/37: coin ranking order changed - Maker and VeChain swapped
# places, each getting freshly updated price/volume figures.
$ws.Range("B36").Value = 'VeChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D36").Value = '0.0165'
$ws.Range("E36").Value = '  +10.21%  '
$ws.Range("B37").Value = 'Maker'
$ws.Range("C37").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D37").Value = '1.121.47'
$ws.Range("E37").Value = '  +3.36%  '

$ws.Range("E38").Value = '  +0.02%  '
$ws.Range("D39").Value = '0.785'
$ws.Range("E39").Value = '  +3.00%  '
$ws.Range("D40").Value = '2.30'
$ws.Range("E40").Value = '  +2.47%  '
$ws.Range("D41").Value = '0.492'
$ws.Range("E41").Value = '  -0.17%  '
$ws.Range("D42").Value = '0.780'
$ws.Range("E42").Value = '  -2.66%  '
$ws.Range("D43").Value = '1.737.97'
$ws.Range("E43").Value = '  +3.46%  '
$ws.Range("D44").Value = '5.14'
$ws.Range("E44").Value = '  +1.92%  '
$ws.Range("D45").Value = '92.91'
$ws.Range("E45").Value = '  +0.14%  '
$ws.Range("D47").Value = '53.66'
$ws.Range("E47").Value = '  +2.55%  '
$ws.Range("D48").Value = '0.0504'
$ws.Range("E48").Value = '  +0.34%  '
$ws.Range("D49").Value = '0.408'
$ws.Range("E49").Value = '  +1.10%  '
$ws.Range("E50").Value = '  -0.22%  '
$ws.Range("D51").Value = '0.0₇0922'
$ws.Range("E51").Value = '  -14.82%  '

# Restore the original "Normal" style so cells don't carry a leftover
# explicit number format that wasn't present in the source file.
$dataRange.Style = "Normal"
